$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from G1 (existing header cell) to H1, then set the value
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("H1").Value = "Save"

$ws.Range("H2").Value = 1
